$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.581.61"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.317.74"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.13"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.62"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.25"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.49"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").Value = "2.729.59"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "56.547.65"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").Value = "2.411.96"
$ws.Range("E17").Value = "  +3.63%  "
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "328.07"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.74"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.23"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("B24").Value = "Kaspa"
$ws.Range("C24").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.164"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.58"
$ws.Range("E25").Value = "  +7.43%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.88"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("D30").Value = "0.0₃0717"
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.09"
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.27"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.882"
$ws.Range("E37").Value = "  -5.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.59"
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.18"
$ws.Range("E40").Value = "  +7.86%  "
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.56"
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "275.56"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.25"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("B49").Value = "Polygon"
$ws.Range("C49").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.380"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0214"
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.04"
$ws.Range("E51").Value = "  +0.89%  "
